$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.465.75"
$ws.Range("E2").Value = "  -5.66%  "

$ws.Range("D3").Value = "1.638.42"
$ws.Range("E3").Value = "  -7.18%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").Value = "'1.007"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").Value = "'305.87"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.30%  "

$ws.Range("D7").Value = "'0.3636"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.02%  "

$ws.Range("D8").Value = "'47.77"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.66%  "

$ws.Range("D9").Value = "'0.3217"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -11.11%  "

$ws.Range("D10").Value = "'1.103"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -10.58%  "

$ws.Range("D11").Value = "'0.06891"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.41%  "

$ws.Range("D12").Value = "'1.011"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.29%  "

$ws.Range("D13").Value = "'5.890"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -9.31%  "

$ws.Range("D14").Value = "'19.16"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -11.67%  "

$ws.Range("D15").Value = "'6.531"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -8.03%  "

$ws.Range("D16").Value = "1.633.08"
$ws.Range("E16").Value = "  -7.33%  "

$ws.Range("D17").Value = "'0.00001043"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -10.17%  "

$ws.Range("D18").Value = "'0.06549"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.82%  "

$ws.Range("D20").Value = "'76.87"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -11.81%  "

$ws.Range("D21").Value = "'15.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -11.33%  "

$ws.Range("D22").Value = "'5.863"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -10.11%  "

$ws.Range("D23").Value = "'11.87"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.29%  "

$ws.Range("D24").Value = "24.396.12"
$ws.Range("E24").Value = "  -5.60%  "

$ws.Range("D25").Value = "'2.417"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("D26").Value = "'2.382"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -19.24%  "

$ws.Range("D27").Value = "'145.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.78%  "

$ws.Range("D28").Value = "'18.74"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -9.53%  "

$ws.Range("D29").Value = "1.821.84"
$ws.Range("E29").Value = "  -7.20%  "

$ws.Range("D30").Value = "'124.44"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -7.26%  "

$ws.Range("D31").Value = "'1.075"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.53%  "

$ws.Range("D32").Value = "'4.084"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.04%  "

$ws.Range("D33").Value = "'5.598"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -22.70%  "

$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").Value = "'0.08384"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.44%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.688"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.66%  "

$ws.Range("D36").Value = "'12.33"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -13.31%  "

$ws.Range("D37").Value = "'5.091"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -11.07%  "

$ws.Range("D38").Value = "'0.06028"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.16%  "

$ws.Range("D39").Value = "'0.02215"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -11.44%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.141"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -13.03%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.192"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.89%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2030"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -10.05%  "

$ws.Range("D43").Value = "'1.007"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.93%  "

$ws.Range("D44").Value = "'0.5855"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -10.92%  "

$ws.Range("D45").Value = "'3.734"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.67%  "

$ws.Range("D46").Value = "'12.45"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -13.62%  "

$ws.Range("E47").Value = "  -12.80%  "

$ws.Range("D48").Value = "'121.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.41%  "

$ws.Range("D49").Value = "'1.917"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -11.73%  "

$ws.Range("D50").Value = "'0.06946"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.43%  "

$ws.Range("D51").Value = "'73.21"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -9.75%  "
